# Apple Mobility Data update for the admissions worksheet.
#
# 1) A data-entry ordering bug swapped the "Sykehuset i Vestfold HF" and
#    "Sykehuset Innlandet HF" rows (out of alphabetical order) for every date
#    block in the existing data. Fix it by swapping the name+value pair back
#    into alphabetical order (Innlandet HF, then i Vestfold HF) for each date.
# 2) Append one more day (2020-04-15, serial 43936) of admissions data for all
#    23 health trusts, extending the sheet from row 875 to row 898.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$pairs = @(
    @(19, 20),
    @(42, 43),
    @(65, 66),
    @(88, 89),
    @(111, 112),
    @(134, 135),
    @(157, 158),
    @(180, 181),
    @(203, 204),
    @(226, 227),
    @(249, 250),
    @(272, 273),
    @(295, 296),
    @(318, 319),
    @(341, 342),
    @(364, 365),
    @(387, 388),
    @(410, 411),
    @(433, 434),
    @(456, 457),
    @(479, 480),
    @(502, 503),
    @(525, 526),
    @(548, 549),
    @(571, 572),
    @(594, 595),
    @(617, 618),
    @(640, 641),
    @(663, 664),
    @(686, 687),
    @(709, 710),
    @(732, 733),
    @(755, 756),
    @(778, 779),
    @(801, 802),
    @(824, 825),
    @(847, 848),
    @(870, 871)
)

foreach ($pair in $pairs) {
    $r1 = $pair[0]
    $r2 = $pair[1]
    $name1 = $ws.Cells.Item($r1, 2).Value2
    $val1  = $ws.Cells.Item($r1, 3).Value2
    $name2 = $ws.Cells.Item($r2, 2).Value2
    $val2  = $ws.Cells.Item($r2, 3).Value2

    $ws.Cells.Item($r1, 2).Value = $name2
    $ws.Cells.Item($r1, 3).Value = $val2
    $ws.Cells.Item($r2, 2).Value = $name1
    $ws.Cells.Item($r2, 3).Value = $val1
}

# New date block: 2020-04-15 (Excel serial 43936)
$newDateSerial = 43936
$newRows = @(
    @(876, "Akershus universitetssykehus HF", 24),
    @(877, "Diakonhjemmet Sykehus", 12),
    @(878, "Finnmarkssykehuset HF", 1),
    @(879, "Haraldsplass Diakonale Sykehus", 4),
    @(880, "Helgelandssykehuset HF", 0),
    @(881, "Helse Bergen HF", 8),
    @(882, "Helse Fonna HF", 4),
    @(883, "Helse Førde HF", 2),
    @(884, "Helse Møre og Romsdal", 5),
    @(885, "Helse Nord-Trøndelag", 4),
    @(886, "Helse Stavanger HF", 6),
    @(887, "Lovisenberg Diakonale Sykehus", 14),
    @(888, "Nordlandssykehuset HF", 6),
    @(889, "Oslo universitetssykehus HF", 32),
    @(890, "Sørlandet sykehus HF", 6),
    @(891, "St. Olavs hospital", 8),
    @(892, "Sunnaas Sykehus HF", 0),
    @(893, "Sykehuset Innlandet HF", 9),
    @(894, "Sykehuset i Vestfold HF", 5),
    @(895, "Sykehuset Østfold HF", 9),
    @(896, "Sykehuset Telemark HF", 3),
    @(897, "Universitetssykehuset Nord-Norge HF", 6),
    @(898, "Vestre Viken HF", 22)
)

foreach ($row in $newRows) {
    $r = $row[0]
    $name = $row[1]
    $val = $row[2]

    $ws.Cells.Item($r, 1).Value = $newDateSerial
    $ws.Cells.Item($r, 1).NumberFormat = "yyyy-mm-dd"
    $ws.Cells.Item($r, 2).Value = $name
    $ws.Cells.Item($r, 3).Value = $val
}

